$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-19 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-20 Monday", 2) | Out-Null
$d.Content.Find.Execute("99×73=7227", $true, $false, $false, $false, $false, $true, 1, $false, "90×27=2430", 2) | Out-Null
$d.Content.Find.Execute("78×84=6552", $true, $false, $false, $false, $false, $true, 1, $false, "69×58=4002", 2) | Out-Null
$d.Content.Find.Execute("22×95=2090", $true, $false, $false, $false, $false, $true, 1, $false, "91×29=2639", 2) | Out-Null
$d.Content.Find.Execute("90×96=8640", $true, $false, $false, $false, $false, $true, 1, $false, "93×23=2139", 2) | Out-Null
$d.Content.Find.Execute("74×35=2590", $true, $false, $false, $false, $false, $true, 1, $false, "39×83=3237", 2) | Out-Null
$d.Content.Find.Execute("56×37=2072", $true, $false, $false, $false, $false, $true, 1, $false, "98×26=2548", 2) | Out-Null
$d.Content.Find.Execute("66×62=4092", $true, $false, $false, $false, $false, $true, 1, $false, "12×70=840", 2) | Out-Null
$d.Content.Find.Execute("49×47=2303", $true, $false, $false, $false, $false, $true, 1, $false, "11×45=495", 2) | Out-Null
$d.Content.Find.Execute("53×17=901", $true, $false, $false, $false, $false, $true, 1, $false, "76×30=2280", 2) | Out-Null
$d.Content.Find.Execute("56×96=5376", $true, $false, $false, $false, $false, $true, 1, $false, "19×26=494", 2) | Out-Null
$d.Content.Find.Execute("43×68=2924", $true, $false, $false, $false, $false, $true, 1, $false, "43×54=2322", 2) | Out-Null
$d.Content.Find.Execute("66×76=5016", $true, $false, $false, $false, $false, $true, 1, $false, "39×66=2574", 2) | Out-Null
$d.Content.Find.Execute("14×22=308", $true, $false, $false, $false, $false, $true, 1, $false, "64×69=4416", 2) | Out-Null
$d.Content.Find.Execute("64×61=3904", $true, $false, $false, $false, $false, $true, 1, $false, "86×38=3268", 2) | Out-Null
$d.Content.Find.Execute("75×33=2475", $true, $false, $false, $false, $false, $true, 1, $false, "31×27=837", 2) | Out-Null
$d.Content.Find.Execute("54×23=1242", $true, $false, $false, $false, $false, $true, 1, $false, "12×20=240", 2) | Out-Null
$d.Content.Find.Execute("68×97=6596", $true, $false, $false, $false, $false, $true, 1, $false, "50×90=4500", 2) | Out-Null
$d.Content.Find.Execute("76×43=3268", $true, $false, $false, $false, $false, $true, 1, $false, "78×94=7332", 2) | Out-Null
$d.Content.Find.Execute("84×52=4368", $true, $false, $false, $false, $false, $true, 1, $false, "47×62=2914", 2) | Out-Null
$d.Content.Find.Execute("80×45=3600", $true, $false, $false, $false, $false, $true, 1, $false, "58×14=812", 2) | Out-Null
$d.Content.Find.Execute("93×31=2883", $true, $false, $false, $false, $false, $true, 1, $false, "61×95=5795", 2) | Out-Null
$d.Content.Find.Execute("97×87=8439", $true, $false, $false, $false, $false, $true, 1, $false, "70×22=1540", 2) | Out-Null
$d.Content.Find.Execute("44×16=704", $true, $false, $false, $false, $false, $true, 1, $false, "20×35=700", 2) | Out-Null
$d.Content.Find.Execute("36×38=1368", $true, $false, $false, $false, $false, $true, 1, $false, "89×53=4717", 2) | Out-Null
$d.Content.Find.Execute("91×32=2912", $true, $false, $false, $false, $false, $true, 1, $false, "99×89=8811", 2) | Out-Null
